# contratos-11-2010.xlsx — fix formatting picked up while scraping floating point numbers.
#
# 1) Three "Razon social"/"Nombre Fantasia" entries used a comma as a name separator,
#    which collides with the CSV/locale comma; switch it to a period.
# 2) The "Importe" column (H) was scraped in es-AR locale formatting
#    ("1.234,56" = thousands "." + decimal ",") instead of plain "1234.56"; rewrite it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Separator fix: comma -> period ---
$ws.Range("E58").Value  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E101").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"

$ws.Range("E60").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F60").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"

$ws.Range("E61").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E102").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# --- 2) "Importe" (column H, rows 2-130): es-AR "1.234,56" -> plain "1234.56" ---
# The values are stored as literal text (General-formatted shared strings), not numbers.
# Mark the range as Text first so Excel keeps the exact digit string we assign instead of
# re-parsing "1037.00" back into a number (and silently dropping the decimals/format).
$ws.Range("H2:H130").NumberFormat = "@"

$ws.Range("H2").Value = "1037.00"
$ws.Range("H3").Value = "1020.00"
$ws.Range("H4").Value = "24999.00"
$ws.Range("H5").Value = "800.00"
$ws.Range("H6").Value = "999.60"
$ws.Range("H7").Value = "143.35"
$ws.Range("H8").Value = "1225.00"
$ws.Range("H9").Value = "84.00"
$ws.Range("H10").Value = "13199.99"
$ws.Range("H11").Value = "52640.00"
$ws.Range("H12").Value = "26798.28"
$ws.Range("H13").Value = "12603.26"
$ws.Range("H14").Value = "6785.90"
$ws.Range("H15").Value = "142.96"
$ws.Range("H16").Value = "5937.76"
$ws.Range("H17").Value = "317.20"
$ws.Range("H18").Value = "256.80"
$ws.Range("H19").Value = "7507.13"
$ws.Range("H20").Value = "300.00"
$ws.Range("H21").Value = "1200.00"
$ws.Range("H22").Value = "860.00"
$ws.Range("H23").Value = "10.80"
$ws.Range("H24").Value = "240.00"
$ws.Range("H25").Value = "5898.03"
$ws.Range("H26").Value = "175.00"
$ws.Range("H27").Value = "179.34"
$ws.Range("H28").Value = "7873.50"
$ws.Range("H29").Value = "120.00"
$ws.Range("H30").Value = "1195.27"
$ws.Range("H31").Value = "15447.22"
$ws.Range("H32").Value = "8665.82"
$ws.Range("H33").Value = "2713.71"
$ws.Range("H34").Value = "22.00"
$ws.Range("H35").Value = "1110.69"
$ws.Range("H36").Value = "2961.40"
$ws.Range("H37").Value = "211.35"
$ws.Range("H38").Value = "2941.65"
$ws.Range("H39").Value = "21.06"
$ws.Range("H40").Value = "4412.72"
$ws.Range("H41").Value = "17.22"
$ws.Range("H42").Value = "367.27"
$ws.Range("H43").Value = "58.00"
$ws.Range("H44").Value = "1800.00"
$ws.Range("H45").Value = "808.00"
$ws.Range("H46").Value = "16.80"
$ws.Range("H47").Value = "3469.55"
$ws.Range("H48").Value = "74.55"
$ws.Range("H49").Value = "120.00"
$ws.Range("H50").Value = "32.00"
$ws.Range("H51").Value = "960.00"
$ws.Range("H52").Value = "1300.00"
$ws.Range("H53").Value = "200.40"
$ws.Range("H54").Value = "12143.20"
$ws.Range("H55").Value = "10225.10"
$ws.Range("H56").Value = "6733.60"
$ws.Range("H57").Value = "1183.00"
$ws.Range("H58").Value = "4032.00"
$ws.Range("H59").Value = "591.00"
$ws.Range("H60").Value = "190.67"
$ws.Range("H61").Value = "3981.00"
$ws.Range("H62").Value = "468.80"
$ws.Range("H63").Value = "154.00"
$ws.Range("H64").Value = "8433.62"
$ws.Range("H65").Value = "4.29"
$ws.Range("H66").Value = "1.54"
$ws.Range("H67").Value = "220445.00"
$ws.Range("H68").Value = "0.68"
$ws.Range("H69").Value = "0.82"
$ws.Range("H70").Value = "5499.18"
$ws.Range("H71").Value = "878.50"
$ws.Range("H72").Value = "1254.00"
$ws.Range("H73").Value = "980.00"
$ws.Range("H74").Value = "481.15"
$ws.Range("H75").Value = "24.00"
$ws.Range("H76").Value = "207.92"
$ws.Range("H77").Value = "11260.00"
$ws.Range("H78").Value = "3120.00"
$ws.Range("H79").Value = "65.80"
$ws.Range("H80").Value = "14320.00"
$ws.Range("H81").Value = "2640.00"
$ws.Range("H82").Value = "865.00"
$ws.Range("H83").Value = "223.49"
$ws.Range("H84").Value = "488.20"
$ws.Range("H85").Value = "182.00"
$ws.Range("H86").Value = "395.00"
$ws.Range("H87").Value = "5000.00"
$ws.Range("H88").Value = "500.00"
$ws.Range("H89").Value = "500.00"
$ws.Range("H90").Value = "2800.00"
$ws.Range("H91").Value = "580.00"
$ws.Range("H92").Value = "500.00"
$ws.Range("H93").Value = "760.00"
$ws.Range("H94").Value = "10382.58"
$ws.Range("H95").Value = "1650.00"
$ws.Range("H96").Value = "1200.00"
$ws.Range("H97").Value = "1500.00"
$ws.Range("H98").Value = "240.00"
$ws.Range("H99").Value = "12080.00"
$ws.Range("H100").Value = "517.00"
$ws.Range("H101").Value = "170.00"
$ws.Range("H102").Value = "40.00"
$ws.Range("H103").Value = "135.00"
$ws.Range("H104").Value = "302.00"
$ws.Range("H105").Value = "6600.01"
$ws.Range("H106").Value = "170.00"
$ws.Range("H107").Value = "102.40"
$ws.Range("H108").Value = "1025.00"
$ws.Range("H109").Value = "1034.78"
$ws.Range("H110").Value = "204.00"
$ws.Range("H111").Value = "11.96"
$ws.Range("H112").Value = "317.30"
$ws.Range("H113").Value = "959.00"
$ws.Range("H114").Value = "3696.00"
$ws.Range("H115").Value = "34.00"
$ws.Range("H116").Value = "216.00"
$ws.Range("H117").Value = "15336.40"
$ws.Range("H118").Value = "17698.70"
$ws.Range("H119").Value = "7650.00"
$ws.Range("H120").Value = "414.96"
$ws.Range("H121").Value = "2545.53"
$ws.Range("H122").Value = "41250.00"
$ws.Range("H123").Value = "546.00"
$ws.Range("H124").Value = "2600.22"
$ws.Range("H125").Value = "831639.22"
$ws.Range("H126").Value = "925.00"
$ws.Range("H127").Value = "389400.00"
$ws.Range("H128").Value = "45000.00"
$ws.Range("H129").Value = "149250.00"
$ws.Range("H130").Value = "1199.23"
